$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header rich-text updates ---
$ws.Range("A8").Characters(21,2).Text = "32"
$ws.Range("C9").Characters(47,8).Text = "8/10/2025"
$ws.Range("C9").Characters(27,9).Text = "8/4/2025"

# --- Cells becoming text placeholders ("0" / "***.*") ---
$ws.Range("C15").Value = "'0"
$ws.Range("G15").Value = "'0"
$ws.Range("H15").Value = "'***.*"
$ws.Range("C18").Value = "'0"
$ws.Range("C27").Value = "'0"
$ws.Range("G27").Value = "'0"
$ws.Range("H27").Value = "'***.*"
$ws.Range("G28").Value = "'0"
$ws.Range("H28").Value = "'***.*"
$ws.Range("C14").Copy()
foreach ($addr in @("C15", "G15", "H15", "C18", "C27", "G27", "H27", "G28", "H28")) { $ws.Range($addr).PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats) }

# --- Numeric value updates ---
$ws.Range("F15").Value = 3
$ws.Range("I15").Value = 10
$ws.Range("K15").Value = -23.076923076923
$ws.Range("L15").Value = 25
$ws.Range("M15").Value = -16.666666666666
$ws.Range("N15").Value = -60
$ws.Range("C16").Value = 1
$ws.Range("C16").NumberFormat = "#,##0"
$ws.Range("D16").Value = 5
$ws.Range("E16").Value = -80
$ws.Range("G16").Value = 8
$ws.Range("H16").Value = -62.5
$ws.Range("I16").Value = 47
$ws.Range("J16").Value = 35
$ws.Range("K16").Value = 34.285714285714
$ws.Range("L16").Value = -20.338983050847
$ws.Range("M16").Value = -76.960784313725
$ws.Range("N16").Value = -92.492012779552
$ws.Range("C17").Value = 3
$ws.Range("D17").Value = 8
$ws.Range("E17").Value = -62.5
$ws.Range("F17").Value = 16
$ws.Range("G17").Value = 17
$ws.Range("H17").Value = -5.882352941176
$ws.Range("I17").Value = 128
$ws.Range("J17").Value = 164
$ws.Range("K17").Value = -21.951219512195
$ws.Range("L17").Value = -15.789473684210
$ws.Range("M17").Value = -31.182795698924
$ws.Range("N17").Value = -47.107438016528
$ws.Range("D18").Value = 1
$ws.Range("E18").Value = -100
$ws.Range("F18").Value = 7
$ws.Range("G18").Value = 12
$ws.Range("H18").Value = -41.666666666666
$ws.Range("J18").Value = 70
$ws.Range("K18").Value = 17.142857142857
$ws.Range("L18").Value = -6.818181818181
$ws.Range("M18").Value = -62.385321100917
$ws.Range("N18").Value = -91.467221644120
$ws.Range("C19").Value = 14
$ws.Range("D19").Value = 5
$ws.Range("E19").Value = 180
$ws.Range("F19").Value = 34
$ws.Range("G19").Value = 26
$ws.Range("H19").Value = 30.769230769230
$ws.Range("I19").Value = 231
$ws.Range("J19").Value = 224
$ws.Range("K19").Value = 3.125
$ws.Range("L19").Value = -10.465116279069
$ws.Range("M19").Value = -12.5
$ws.Range("N19").Value = -35.833333333333
$ws.Range("C20").Value = 8
$ws.Range("D20").Value = 9
$ws.Range("E20").Value = -11.111111111111
$ws.Range("F20").Value = 17
$ws.Range("G20").Value = 29
$ws.Range("H20").Value = -41.379310344827
$ws.Range("I20").Value = 117
$ws.Range("J20").Value = 151
$ws.Range("K20").Value = -22.516556291390
$ws.Range("L20").Value = 23.157894736842
$ws.Range("M20").Value = -49.568965517241
$ws.Range("N20").Value = -94.309338521400
$ws.Range("C21").Value = 26
$ws.Range("D21").Value = 28
$ws.Range("E21").Value = -7.142857142857
$ws.Range("F21").Value = 80
$ws.Range("G21").Value = 92
$ws.Range("H21").Value = -13.043478260869
$ws.Range("I21").Value = 616
$ws.Range("J21").Value = 657
$ws.Range("K21").Value = -6.240487062404
$ws.Range("L21").Value = -6.807866868381
$ws.Range("M21").Value = -45.244444444444
$ws.Range("N21").Value = -85.614198972442
$ws.Range("C24").Value = 49
$ws.Range("D24").Value = 11
$ws.Range("E24").Value = 345.454545454545
$ws.Range("F24").Value = 81
$ws.Range("G24").Value = 48
$ws.Range("H24").Value = 68.75
$ws.Range("I24").Value = 391
$ws.Range("J24").Value = 423
$ws.Range("K24").Value = -7.565011820330
$ws.Range("L24").Value = -18.541666666666
$ws.Range("M24").Value = -22.111553784860
$ws.Range("C25").Value = 1
$ws.Range("D25").Value = 7
$ws.Range("E25").Value = -85.714285714285
$ws.Range("F25").Value = 12
$ws.Range("G25").Value = 14
$ws.Range("H25").Value = -14.285714285714
$ws.Range("I25").Value = 92
$ws.Range("J25").Value = 109
$ws.Range("K25").Value = -15.596330275229
$ws.Range("L25").Value = 6.976744186046
$ws.Range("C26").Value = 6
$ws.Range("D26").Value = 9
$ws.Range("E26").Value = -33.333333333333
$ws.Range("F26").Value = 25
$ws.Range("G26").Value = 30
$ws.Range("H26").Value = -16.666666666666
$ws.Range("I26").Value = 257
$ws.Range("J26").Value = 260
$ws.Range("K26").Value = -1.153846153846
$ws.Range("L26").Value = 21.800947867298
$ws.Range("M26").Value = -35.101010101010
$ws.Range("F27").Value = 3
$ws.Range("I27").Value = 10
$ws.Range("K27").Value = -44.444444444444
$ws.Range("L27").Value = -16.666666666666
$ws.Range("C28").Value = 3
$ws.Range("C28").NumberFormat = "#,##0"
$ws.Range("F28").Value = 3
$ws.Range("I28").Value = 23
$ws.Range("K28").Value = -4.166666666666
$ws.Range("L28").Value = 64.285714285714
$ws.Range("L33").Value = -50
